# Update timestamps in the "data" sheet (column F) to reflect the refreshed
# panel query time.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = "2021-10-05 14:35:05.856244"
$ws1.Range("F3").Value = "2021-10-05 14:35:05.856252"
$ws1.Range("F4").Value = "2021-10-05 14:35:05.856256"
$ws1.Range("F5").Value = "2021-10-05 14:35:05.856259"
$ws1.Range("F6").Value = "2021-10-05 14:35:05.856261"
$ws1.Range("F7").Value = "2021-10-05 14:35:05.856264"
$ws1.Range("F8").Value = "2021-10-05 14:35:05.856266"
$ws1.Range("F9").Value = "2021-10-05 14:35:05.856269"
$ws1.Range("F10").Value = "2021-10-05 14:35:05.856272"
$ws1.Range("F11").Value = "2021-10-05 14:35:05.856274"
$ws1.Range("F12").Value = "2021-10-05 14:35:05.856277"
$ws1.Range("F13").Value = "2021-10-05 14:35:05.856279"
$ws1.Range("F14").Value = "2021-10-05 14:35:05.856282"
$ws1.Range("F15").Value = "2021-10-05 14:35:05.856284"
$ws1.Range("F16").Value = "2021-10-05 14:35:05.856287"
$ws1.Range("F17").Value = "2021-10-05 14:35:05.856289"
$ws1.Range("F18").Value = "2021-10-05 14:35:05.856292"
$ws1.Range("F19").Value = "2021-10-05 14:35:05.856294"
$ws1.Range("F20").Value = "2021-10-05 14:35:05.856297"
$ws1.Range("F21").Value = "2021-10-05 14:35:05.856299"

# Add a new "metadata" worksheet, placed after the "data" sheet, describing
# the panel query that produced the "data" sheet contents.
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "metadata"

# Reuse the bold/bordered/centered header style already used for the "data"
# sheet's header row (style index 1) instead of minting a new one.
$ws1.Range("E1:F1").Copy()
$newSheet.Range("F1:G1").PasteSpecial(-4122)
$ws1.Range("B1:F1").Copy()
$newSheet.Range("B1:F1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Header row.
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row describing the query for this panel.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Ocular and Oculocutaneous Albinism"
$newSheet.Range("C2").Value = 37
$newSheet.Range("D2").Formula = '=TEXT(1,"0.0")'
$newSheet.Range("D2").Copy()
$newSheet.Range("D2").PasteSpecial(-4163)
$newSheet.Range("E2").Value = "2021-06-04T10:07:27.580689Z"
$newSheet.Range("F2").Value = "2021-10-05 14:35:05.852338"
$newSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/37/?format=json"
